$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row (only within the used columns A:T) before row 22,
# shifting rows 22:69 down to 23:70. Using a bounded range (rather than a
# full-row insert) avoids touching unused columns / bloating the sheet.
$ws.Range("A22:T22").Insert(-4121)  # xlShiftDown

# Copy formatting from the (now shifted) row below (row 23, originally row 22)
# so the new row 22 keeps the same styles (e.g. date format in column D).
$ws.Range("A23:T23").Copy()
$ws.Range("A22:T22").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 22 with the data for the new entry
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44720
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100108
$ws.Range("H22").Value = "Tropicales y subtropicales"
$ws.Range("I22").Value = 100108007
$ws.Range("J22").Value = "Coco"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 28000
$ws.Range("O22").Value = 28000
$ws.Range("P22").Value = 28000
$ws.Range("Q22").Value = "$/malla 20 unidades"
$ws.Range("R22").Value = "Perú"
$ws.Range("S22").Value = 1400
$ws.Range("T22").Value = 20
